# "fixed R49 from 510k into 510"
#
# The BOM had R49 duplicated across two rows:
#   - a row with Value "510k" whose Parts column only listed "R49"
#   - a row with Value "510" (a plain number) whose Parts column listed
#     "R9, R10, R37"
# R49 actually belongs with the 510-ohm resistors, not the 510k ones, so
# the "510k" row is removed entirely and R49 is appended to the Parts list
# of the "510" row (bumping its Qty from 3 to 4).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the stray "510k" / "R49" row (row 20).
$ws.Rows("20:20").Delete()

# The "510" row (old row 23) is now row 22 after the shift above.
# Update its quantity and merge R49 into its Parts list.
$ws.Range("A22").Value = 4
$ws.Range("E22").Value = "R9, R10, R37, R49"

# Match the saved selection/cursor position recorded in the workbook.
$ws.Range("A23").Select() | Out-Null
